# Generate Report for Handback
#
# Two logical "handback" records are refreshed in this run:
#   - the existing row for file "2175c337-...md" is replaced in-place by a
#     newer run against a renamed/regenerated source file
#     "ac55d992-1edd-410e-bdfd-f7bebed1963e.md"
#   - a brand new second file "c9c70c95-c479-43fd-a4d2-48270d45b7c0.md" is
#     handed back for the first time, appended as a new row on every sheet.
#
# This script reproduces both changes across the Overview / zh-cn / de-de
# worksheets, keeping their backing tables in sync.

$wb = $excel.ActiveWorkbook

$file1Old = "2175c337-6dc0-4918-aa9a-89e15c8d7752"
$file1New = "ac55d992-1edd-410e-bdfd-f7bebed1963e"
$file2New = "c9c70c95-c479-43fd-a4d2-48270d45b7c0"

$hoHashZh = "345e5284197f96f7f8137abb0f1f643e2aedff4e"
$hoHashDe = "345e5284197f96f7f8137abb0f1f643e2aedff4e"
$hoHashZh2 = "5ce180b0445c7c4ad289ac1955311cb07e44f947"
$hoHashDe2 = "5ce180b0445c7c4ad289ac1955311cb07e44f947"

$statusText = "Handed back: in sync with en-US"

$overviewDate = "2016-08-24 07:02:42"
$zhHoDate = "2016-08-24 07:02:36"
$zhHbDate1 = "2016-08-24 07:03:05"
$zhHbDate2 = "2016-08-24 07:03:05"
$deHbDate1 = "2016-08-24 07:03:17"
$deHbDate2 = "2016-08-24 07:03:17"

$repoBase = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/81bd80e356ab52cf8c73078ba640afbc66dfebf1/e2e"
$repoZh = "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/06d9b6bb345add54db486e4ddf7d2f31105938a2/e2e"
$repoDe = "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/d918679705724975eb5f0ce682b7eebebf6d52be/e2e"

# ----------------------------------------------------------------------
# Overview sheet
# ----------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")

# Row 2: refresh the existing record (2175c337... -> ac55d992...)
$ws.Range("A2").Value = "$file1New.md"
$ws.Range("B2").Value = "e2e\$file1New.md"
$ws.Hyperlinks.Add($ws.Range("B2"), "$repoBase/$file1New.md", "", "", "e2e\$file1New.md") | Out-Null
$ws.Range("C2").Value = ".md"
$ws.Range("E2").Value = $statusText
$ws.Range("F2").Value = $statusText
$ws.Range("G2").Value = $overviewDate

# Row 3: brand new record (c9c70c95...)
$ws.Range("A3").Value = "$file2New.md"
$ws.Range("B3").Value = "e2e\$file2New.md"
$ws.Hyperlinks.Add($ws.Range("B3"), "$repoBase/$file2New.md", "", "", "e2e\$file2New.md") | Out-Null
$ws.Range("C3").Value = ".md"
$ws.Range("E3").Value = $statusText
$ws.Range("F3").Value = $statusText
$ws.Range("G3").Value = $overviewDate

$loOverview = $ws.ListObjects.Item(1)
$loOverview.Resize($ws.Range("A1:G3"))

# ----------------------------------------------------------------------
# zh-cn sheet
# ----------------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")

# Row 2: refresh the existing record
$ws.Range("A2").Value = "$file1New.md"
$ws.Hyperlinks.Add($ws.Range("A2"), "$repoBase/$file1New.md", "", "", "$file1New.md") | Out-Null
$ws.Range("B2").Value = ".md"
$ws.Range("C2").Value = $statusText
$ws.Range("D2").Value = "e2e"
$ws.Range("E2").Value = "ht"
$ws.Range("F2").Value = "False"
$ws.Range("G2").Value = "$file1New.$hoHashZh.zh-cn.xlf"
$ws.Range("H2").Value = $zhHoDate
$ws.Range("I2").Value = "$file1New.md"
$ws.Hyperlinks.Add($ws.Range("I2"), "$repoZh/$file1New.md", "", "", "$file1New.md") | Out-Null
$ws.Range("J2").Value = "$file1New.$hoHashZh.zh-cn.xlf"
$ws.Range("K2").Value = $zhHbDate1
$ws.Range("L2").Value = ""
$ws.Range("M2").Value = "True"
$ws.Range("N2").Value = ""
$ws.Range("O2").Value = "False"
$ws.Range("P2").Value = ""

# Row 3: brand new record
$ws.Range("A3").Value = "$file2New.md"
$ws.Hyperlinks.Add($ws.Range("A3"), "$repoBase/$file2New.md", "", "", "$file2New.md") | Out-Null
$ws.Range("B3").Value = ".md"
$ws.Range("C3").Value = $statusText
$ws.Range("D3").Value = "e2e"
$ws.Range("E3").Value = "ht"
$ws.Range("F3").Value = "True"
$ws.Range("G3").Value = "$file2New.$hoHashZh2.zh-cn.xlf"
$ws.Range("H3").Value = $zhHoDate
$ws.Range("I3").Value = "$file2New.md"
$ws.Hyperlinks.Add($ws.Range("I3"), "$repoZh/$file2New.md", "", "", "$file2New.md") | Out-Null
$ws.Range("J3").Value = "$file2New.$hoHashZh2.zh-cn.xlf"
$ws.Range("K3").Value = $zhHbDate2
$ws.Range("L3").Value = ""
$ws.Range("M3").Value = "True"
$ws.Range("N3").Value = ""
$ws.Range("O3").Value = "False"
$ws.Range("P3").Value = ""

$loZh = $ws.ListObjects.Item(1)
$loZh.Resize($ws.Range("A1:P3"))

# ----------------------------------------------------------------------
# de-de sheet
# ----------------------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")

# Row 2: refresh the existing record
$ws.Range("A2").Value = "$file1New.md"
$ws.Hyperlinks.Add($ws.Range("A2"), "$repoBase/$file1New.md", "", "", "$file1New.md") | Out-Null
$ws.Range("B2").Value = ".md"
$ws.Range("C2").Value = $statusText
$ws.Range("D2").Value = "e2e"
$ws.Range("E2").Value = "ht"
$ws.Range("F2").Value = "False"
$ws.Range("G2").Value = "$file1New.$hoHashDe.de-de.xlf"
$ws.Range("H2").Value = $overviewDate
$ws.Range("I2").Value = "$file1New.md"
$ws.Hyperlinks.Add($ws.Range("I2"), "$repoDe/$file1New.md", "", "", "$file1New.md") | Out-Null
$ws.Range("J2").Value = "$file1New.$hoHashDe.de-de.xlf"
$ws.Range("K2").Value = $deHbDate1
$ws.Range("L2").Value = ""
$ws.Range("M2").Value = "True"
$ws.Range("N2").Value = ""
$ws.Range("O2").Value = "False"
$ws.Range("P2").Value = ""

# Row 3: brand new record
$ws.Range("A3").Value = "$file2New.md"
$ws.Hyperlinks.Add($ws.Range("A3"), "$repoBase/$file2New.md", "", "", "$file2New.md") | Out-Null
$ws.Range("B3").Value = ".md"
$ws.Range("C3").Value = $statusText
$ws.Range("D3").Value = "e2e"
$ws.Range("E3").Value = "ht"
$ws.Range("F3").Value = "True"
$ws.Range("G3").Value = "$file2New.$hoHashDe2.de-de.xlf"
$ws.Range("H3").Value = $overviewDate
$ws.Range("I3").Value = "$file2New.md"
$ws.Hyperlinks.Add($ws.Range("I3"), "$repoDe/$file2New.md", "", "", "$file2New.md") | Out-Null
$ws.Range("J3").Value = "$file2New.$hoHashDe2.de-de.xlf"
$ws.Range("K3").Value = $deHbDate2
$ws.Range("L3").Value = ""
$ws.Range("M3").Value = "True"
$ws.Range("N3").Value = ""
$ws.Range("O3").Value = "False"
$ws.Range("P3").Value = ""

$loDe = $ws.ListObjects.Item(1)
$loDe.Resize($ws.Range("A1:P3"))

Write-Output "Handback report regenerated: $file1New.md refreshed, $file2New.md added."
